$d = $word.ActiveDocument

# Locate the target paragraph: "The WediOnAccess component is invoked when the user
# access an URL. Also this component attaches global event listeners for interacting
# with the rest components and buttons from browser."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*WediOnAccess component is invoked*") {
        $target = $p
        break
    }
}

# --- 1. Split " component " into " " + "component " (both keep the same
#        Comic Sans MS / 28pt formatting the original run already had). A
#        genuine no-op font assignment does not force a run split, so we
#        nudge the name away and back to guarantee the split happens.
$r = $target.Range.Duplicate
$r.Find.Execute("component ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, 0) | Out-Null
$r.Font.Name = "Arial"
$r.Font.Name = "Comic Sans MS"

# --- 2. "is invoked " run: add Comic Sans MS font and shrink size 29 -> 28.
$r = $target.Range.Duplicate
$r.Find.Execute("is invoked ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, 0) | Out-Null
$r.Font.Name = "Comic Sans MS"
$r.Font.Size = 14

# --- 3. Move the lone "_GoBack" bookmark from the end of the document (after
#        the inline image) to just before "when the user access ".
$r = $target.Range.Duplicate
$r.Find.Execute("when the user access ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, 0) | Out-Null
$bmRange = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- 4. "when the user access " run: Comic Sans MS + 28pt.
$r = $target.Range.Duplicate
$r.Find.Execute("when the user access ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, 0) | Out-Null
$r.Font.Name = "Comic Sans MS"
$r.Font.Size = 14

# --- 5. "an" run.
$r = $target.Range.Duplicate
$r.Find.Execute("an", $false, $true, $false, $false, $false, $true, 1, $false, "", 0, 0) | Out-Null
$r.Font.Name = "Comic Sans MS"
$r.Font.Size = 14

# --- 6. " URL. Also this component " run.
$r = $target.Range.Duplicate
$r.Find.Execute(" URL. Also this component ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, 0) | Out-Null
$r.Font.Name = "Comic Sans MS"
$r.Font.Size = 14

# --- 7. "attaches global event listeners for i" run.
$r = $target.Range.Duplicate
$r.Find.Execute("attaches global event listeners for i", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, 0) | Out-Null
$r.Font.Name = "Comic Sans MS"
$r.Font.Size = 14

# --- 8. "nteracting with the rest components and buttons from browser" run.
$r = $target.Range.Duplicate
$r.Find.Execute("nteracting with the rest components and buttons from browser", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, 0) | Out-Null
$r.Font.Name = "Comic Sans MS"
$r.Font.Size = 14

# --- 9. Final "." run.
$r = $target.Range.Duplicate
$r.Find.Execute("browser.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, 0) | Out-Null
$r2 = $d.Range($r.End - 1, $r.End)
$r2.Font.Name = "Comic Sans MS"
$r2.Font.Size = 14

Write-Output "done"
